# Regenerate save_data: column G ("K") was recomputed using the Strike
# count (K) in place of the previous Strike# values. Rewrite column G for
# every data row (rows 2-63) with the recalculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(1,1,1,2,2,1,0,3,1,1,2,2,1,2,1,1,0,1,1,1,0,1,0,2,1,1,1,0,2,1,1,0,0,0,2,0,1,0,1,1,1,3,1,0,1,0,1,1,3,0,1,0,0,0,1,0,1,1,0,0,2,0)

$row = 2
foreach ($v in $kValues) {
    $ws.Cells.Item($row, 7).Value = $v
    $row++
}
